$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as literal TEXT (not auto-parsed as a number) without
# leaving the cell's number-format permanently altered. We flip the format to
# "@" just long enough to assign the literal, then paste-special the FORMAT
# (not the value) from a pristine, never-reformatted cell back on top so the
# on-disk style index returns to the sheet's default (matches how the other
# quarters store numeric-looking figures as plain inline/shared strings).
# ---------------------------------------------------------------------------
function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet, positioned right after "2021-Q4" and right
#    before "总计" — duplicate "2021-Q4" so the layout/styles match the other
#    quarterly sheets exactly, then overwrite it with the new quarter's data.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newQ = $wb.Worksheets.Item($q4.Index + 1)
$newQ.Name = "2022-Q1"

Set-TextValue $newQ.Range("B2") "008763"
Set-TextValue $newQ.Range("C2") "天弘越南市场股票（QDII）A"
Set-TextValue $newQ.Range("D2") "37.53"
Set-TextValue $newQ.Range("E2") "92.10"
Set-TextValue $newQ.Range("F2") "8.48"
Set-TextValue $newQ.Range("G2") "3.1825"
$newQ.Range("H2").Value = 1

Set-TextValue $newQ.Range("B3") "008764"
Set-TextValue $newQ.Range("C3") "天弘越南市场股票（QDII）C"
Set-TextValue $newQ.Range("D3") "14.26"
Set-TextValue $newQ.Range("E3") "92.10"
Set-TextValue $newQ.Range("F3") "8.48"
Set-TextValue $newQ.Range("G3") "1.2092"
$newQ.Range("H3").Value = 1

# Restore the default (un-styled) look for the cells we just touched by
# stamping the format of an untouched numeric cell (H2, style 0) over them.
$newQ.Range("H2").Copy()
$newQ.Range("B2:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update "总计": prepend a new row for 2022-Q1 (2 holdings, 4.39 亿元),
#    pushing the existing quarters down and renumbering the index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 4.39

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Re-apply the index column's usual style (border/bold) to the new row's A2.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Restore the original active sheet/selection (sheet-copy operations above
#    shift Excel's focus onto the newly created sheet as a side effect).
# ---------------------------------------------------------------------------
$first = $wb.Worksheets.Item("2021-Q2")
$first.Select() | Out-Null
$first.Range("A1").Select() | Out-Null
